$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row heights for rows 1-7 (content/layout reflowed) ---
$ws.Rows.Item(1).RowHeight = 63
$ws.Rows.Item(2).RowHeight = 63
$ws.Rows.Item(3).RowHeight = 47.25
$ws.Rows.Item(4).RowHeight = 63
$ws.Rows.Item(5).RowHeight = 47.25
$ws.Rows.Item(6).RowHeight = 94.5
$ws.Rows.Item(7).RowHeight = 94.5

# --- Insert 40 new rows starting at row 15 ---
$ws.Rows.Item(15).Resize(40).Insert()

# --- Rows 15-34: term/meaning pairs (entered A/B together) ---
$ws.Range("A15").Value = '上げる'
$ws.Range("B15").Value = 'ageru – Làm…xong'
$ws.Range("A16").Value = '合う'
$ws.Range("B16").Value = 'au – Làm điều gì đó cùng nhau'
$ws.Range("A17").Value = 'ばよかった'
$ws.Range("B17").Value = 'ba yokatta – Giá mà… thì tốt rồi'
$ws.Range("A18").Value = '〜ほど'
$ws.Range("B18").Value = 'ba~hodo – Càng…càng'
$ws.Range("A19").Value = '〜のに'
$ws.Range("B19").Value = 'ba~noni – Thế mà/giá mà'
$ws.Range("A20").Value = 'ばかり'
$ws.Range("B20").Value = 'bakari – Toàn…/chỉ…/lúc nào cũng'
$ws.Range("A21").Value = 'ばかりか〜も'
$ws.Range("B21").Value = 'bakarika~mo – Không chỉ có… mà còn'
$ws.Range("A22").Value = 'べきだ'
$ws.Range("B22").Value = 'beki da – Phải/nên làm gì'
$ws.Range("A23").Value = '別に〜ない'
$ws.Range("B23").Value = 'betsuni ni~nai – Không thực sự/không hẳn'
$ws.Range("A24").Value = 'ぶりに'
$ws.Range("B24").Value = 'buri ni – Sau (mới lại)…'
$ws.Range("A25").Value = 'ちゃった'
$ws.Range("B25").Value = 'chatta – Làm gì đó mất rồi'
$ws.Range("A26").Value = 'だけ'
$ws.Range("B26").Value = 'dake – đến mức tối đa có thể…/được chứng nào hay chứng đó'
$ws.Range("A27").Value = 'だけでなく'
$ws.Range("B27").Value = 'dake de naku – Không chỉ…mà còn'
$ws.Range("A28").Value = 'だけしか'
$ws.Range("B28").Value = 'dake shika – Chỉ… mà thôi'
$ws.Range("A29").Value = 'だけど'
$ws.Range("B29").Value = 'dakedo – Nhưng mà'
$ws.Range("A30").Value = 'だから'
$ws.Range("B30").Value = 'desu kara – Vì vậy'
$ws.Range("A31").Value = 'どんなに〜ても'
$ws.Range("B31").Value = 'donnani~temo – Cho dù có như thế nào/bao nhiêu đi nữa'
$ws.Range("A32").Value = 'どうしても'
$ws.Range("B32").Value = 'doushitemo – Bằng bất cứ giá nào/dù thế nào đi chăng nữa'
$ws.Range("A33").Value = 'ふりをする'
$ws.Range("B33").Value = 'furi wo suru – Giả vờ làm gì/Tỏ ra như thể là…'
$ws.Range("A34").Value = 'ふと'
$ws.Range("B34").Value = 'futo – Đột nhiên/bất ngờ/chợt'

# --- Rows 35-54: grammar column A (all terms first) ---
$ws.Range("A35").Value = 'がる'
$ws.Range("A36").Value = 'ごらん'
$ws.Range("A37").Value = 'ごとに'
$ws.Range("A38").Value = 'はずだ'
$ws.Range("A39").Value = 'ほど'
$ws.Range("A40").Value = 'ほど～ない'
$ws.Range("A41").Value = '一度に'
$ws.Range("A42").Value = 'いくら～ても'
$ws.Range("A43").Value = '一体'
$ws.Range("A44").Value = 'か何か'
$ws.Range("A45").Value = 'かける'
$ws.Range("A46").Value = 'かなあ'
$ws.Range("A47").Value = '必ずしも～とは限らない'
$ws.Range("A48").Value = 'から～にかけて'
$ws.Range("A49").Value = '代わりに'
$ws.Range("A50").Value = '結局'
$ws.Range("A51").Value = '決して～ない'
$ws.Range("A52").Value = 'きり'
$ws.Range("A53").Value = '切る/切れる/切れない'
$ws.Range("A54").Value = 'つけ'

# --- Rows 35-54: meaning column B (Japanese<br>Vietnamese, filled afterward) ---
$ws.Range("B35").Value = 'がる<br>Ý muốn (của người khác)'
$ws.Range("B36").Value = 'ごらん<br>Hãy làm/thử xem...'
$ws.Range("B37").Value = 'ごとに<br>Mỗi/cứ mỗi/cứ...lại...'
$ws.Range("B38").Value = 'はずだ<br>Chắc chắn là...'
$ws.Range("B39").Value = 'ほど<br>Thường...(hơn)'
$ws.Range("B40").Value = 'ほど～ない<br>Không bằng như.../không tới mức như...'
$ws.Range("B41").Value = 'いちどに<br>Cùng một lúc/tất cả trong một ...'
$ws.Range("B42").Value = 'いくら～ても<br>Bất kể thế nào/dù thế nào đi chăng nữa ...'
$ws.Range("B43").Value = 'いったい<br>Rốt cuộc/ không biết là/ vậy thì'
$ws.Range("B44").Value = 'か なに か<br>Hay gì đó'
$ws.Range("B45").Value = 'かける<br>Chưa xong/dở dang'
$ws.Range("B46").Value = 'かなあ<br>Hay sao/mong sao/có ... không đây'
$ws.Range("B47").Value = 'かならずしも～とはかぎらない<br>Không nhất thiết là ...'
$ws.Range("B48").Value = 'から～にかけて<br>Từ ... đến'
$ws.Range("B49").Value = 'かわりに<br>Thay vì/đổi lại/thay cho'
$ws.Range("B50").Value = 'けっきょく<br>Sau tất cả/ cuối cùng'
$ws.Range("B51").Value = 'けっして～ない<br>Nhất định không/tuyệt đối không'
$ws.Range("B52").Value = 'きり<br>Chỉ có/có'
$ws.Range("B53").Value = 'きる / きれる / きれない<br>Hết/ không hết'
$ws.Range("B54").Value = 'つけ<br>Có phải ... đúng không/có phải là'

# --- Row heights within the new block that differ from the default ---
$ws.Rows.Item(26).RowHeight = 31.5
$ws.Rows.Item(31).RowHeight = 31.5
$ws.Rows.Item(32).RowHeight = 31.5
$ws.Rows.Item(33).RowHeight = 31.5
$ws.Rows.Item(40).RowHeight = 31.5
$ws.Rows.Item(42).RowHeight = 31.5
$ws.Rows.Item(46).RowHeight = 31.5
$ws.Rows.Item(47).RowHeight = 31.5
$ws.Rows.Item(51).RowHeight = 31.5

# --- Apply wrap + vertical-center style to the whole new block (A15:B54) ---
$rngMain = $ws.Range("A15:B54")
$rngMain.WrapText = $true
$rngMain.VerticalAlignment = -4108

# --- A15 keeps the existing yellow highlight (same fill as A7) ---
$ws.Range("A15").Interior.Color = 65535

# --- Column C: width + matching wrap/vcenter style for rows 35-54 ---
$ws.Columns.Item(3).ColumnWidth = 37.21875
$rngC = $ws.Range("C35:C54")
$rngC.WrapText = $true
$rngC.VerticalAlignment = -4108

# --- Sheet view: scroll position + active selection ---
$ws.Activate()
$ws.Range("K31").Select()
try { $excel.ActiveWindow.ScrollRow = 28 } catch {}
try { $excel.ActiveWindow.TopLeftCell = $ws.Range("A28") } catch {}

# --- Conditional formatting: keep duplicate-check off the newly inserted rows ---
try {
    $cfItem = $ws.Range("A1:A1048576").FormatConditions.Item(1)
    $cfItem.ModifyAppliesToRange($ws.Range("A1:A14,A55:A1048576"))
} catch {}

Write-Output "edit complete"